$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the plain/default style from an existing unstyled cell so we can
# re-apply it after writing text values (PasteSpecial keeps them as text
# instead of Excel auto-converting numeric-looking strings like "$ 0" or
# "4510.111" into numbers).
$ws.Cells.Item(2, 1).Copy() | Out-Null

$rowsData = @(
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'DUS UCM', "'" + '4510.111', "'" + 'UCM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'TAIL LIGHTS-REAR LICENSE PLATE', "'" + '4513.05', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'DUS UCM', "'" + '4510.111', "'" + 'UCM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'TAIL LIGHTS-REAR LICENSE PLATE', "'" + '4513.05', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'DUS UCM', "'" + '4510.111', "'" + 'UCM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'TAIL LIGHTS-REAR LICENSE PLATE', "'" + '4513.05', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'DUS UCM', "'" + '4510.111', "'" + 'UCM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21TRD09386', "'" + 'Hemmeter', "'" + 'TAIL LIGHTS-REAR LICENSE PLATE', "'" + '4513.05', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'No Contest', "'" + 'Guilty', "'" + '$ 0', "'" + '$ 0')
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'Not Guilty', $null, $null, $null)
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'Not Guilty', $null, $null, $null)
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'Not Guilty', $null, $null, $null)
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'Not Guilty', $null, $null, $null)
    ,@("'" + '21CRB01291', "'" + 'Hemmeter', "'" + 'PERMISSION REQ''D TO USE LICENSED DOCK', "'" + '1501:46-12-04', "'" + 'MM', "'" + 'Not Guilty', $null, $null, $null)
)

$startRow = 857
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $rowsData[$i]
    for ($c = 0; $c -lt 9; $c++) {
        $val = $rowVals[$c]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($r, $c + 1)
            $cell.Value = $val
            $cell.PasteSpecial(-4122) | Out-Null
        }
    }
}

$excel.CutCopyMode = $false
$ws.Range("B2").Select() | Out-Null